# Raw and Clean Data From SSA for July 14th
#
# The historic log sheet ("out_vars") gets one new daily row (45) appended
# for 2020-07-14, with its five metric columns (Confirmados, Negativos,
# Sospechosos, Defunciones, Porcentaje hospitalizados).
#
# While here, also clear the stray "Text" number format that had been
# applied to A35:A44 (the date cells for 07-04..07-13) so the whole date
# column is back to a uniform, unformatted style - matching the rest of
# column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the leftover Text number format on the existing date cells ---
for ($r = 35; $r -le 44; $r++) {
    $ws.Cells.Item($r, 1).Style = "Normal"
}

# --- Append the new row for 2020-07-14 ---
$newRow = 45

# Force the date to be stored as text (matches every other date cell in
# column A, which are shared strings, not real dates) and then drop the
# format back to the default style once the text value has been set.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2020-07-14"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 311486
$ws.Cells.Item($newRow, 3).Value = 363930
$ws.Cells.Item($newRow, 4).Value = 80721
$ws.Cells.Item($newRow, 5).Value = 36327
$ws.Cells.Item($newRow, 6).Value = 29.12
